$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'" + "42.702.78"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Formula = "'" + "2.302.61"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Formula = "'" + "316.76"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Formula = "'" + "103.80"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Formula = "'" + "0.628"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Formula = "'" + "39.81"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Formula = "'" + "0.0905"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Formula = "'" + "8.50"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Formula = "'" + "0.999"
$ws.Range("E14").Value = "  +3.34%  "
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Formula = "'" + "2.302.07"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Formula = "'" + "42.623.83"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Formula = "'" + "7.58"
$ws.Range("E19").Value = "  +2.67%  "
$ws.Range("D20").Formula = "'" + "13.84"
$ws.Range("E20").Value = "  +28.15%  "
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").Formula = "'" + "74.07"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Formula = "'" + "266.45"
$ws.Range("E24").Value = "  -4.51%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Formula = "'" + "2.35"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Formula = "'" + "22.65"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Formula = "'" + "6.63"
$ws.Range("E30").Value = "  +13.42%  "
$ws.Range("D31").Formula = "'" + "37.66"
$ws.Range("E31").Value = "  +4.29%  "
$ws.Range("D32").Formula = "'" + "165.81"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").Formula = "'" + "0.131"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Formula = "'" + "4.58"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  +11.34%  "
$ws.Range("D42").Formula = "'" + "70.66"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").Formula = "'" + "95.24"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Formula = "'" + "12.49"
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("D47").Formula = "'" + "117.98"
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("D48").Formula = "'" + "80.17"
$ws.Range("E48").Value = "  +3.66%  "
$ws.Range("D49").Formula = "'" + "1.657.85"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("D50").Formula = "'" + "8.90"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -0.06%  "
